# Generate Report for Archive
# Swap the row-4 / row-5 data (for 8de05a06-... and e2c142c8-... entries) on all
# three worksheets, and refresh the status of the e2c142c8 entry (now in row 4)
# from "Ready for handoff" to "In Translation". Hyperlink target r:id's stay
# fixed per row/cell position, only their displayed text needs to follow the
# data that now lives in that row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Overview" sheet (columns A:G)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "e2c142c8-972f-4f3e-9189-484965d728c7.md"
$wsOverview.Range("B4").Value = "e2e\e2c142c8-972f-4f3e-9189-484965d728c7.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "In Translation"
$wsOverview.Range("F4").Value = "In Translation"
$wsOverview.Range("G4").Value = "2016-10-19 16:09:25"

$wsOverview.Range("A5").Value = "8de05a06-e841-430d-ad60-ba25ede17482.md"
$wsOverview.Range("B5").Value = "e2e\8de05a06-e841-430d-ad60-ba25ede17482.md"
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-10-19 16:10:00"

foreach ($h in $wsOverview.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$4') {
        $h.TextToDisplay = "e2e\e2c142c8-972f-4f3e-9189-484965d728c7.md"
    } elseif ($addr -eq '$B$5') {
        $h.TextToDisplay = "e2e\8de05a06-e841-430d-ad60-ba25ede17482.md"
    }
}

# ---------------------------------------------------------------------------
# "zh-cn" sheet (columns A:P)
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A4").Value = "e2c142c8-972f-4f3e-9189-484965d728c7.md"
$wsZhCn.Range("C4").Value = "In Translation"
$wsZhCn.Range("G4").Value = "e2c142c8-972f-4f3e-9189-484965d728c7.8bc76d304de52b61f9888771bb68e840895cbd08.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-10-19 16:09:13"

$wsZhCn.Range("A5").Value = "8de05a06-e841-430d-ad60-ba25ede17482.md"
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("G5").Value = "8de05a06-e841-430d-ad60-ba25ede17482.205adbe7658bb18bc81004be3c10cc68a87a9472.zh-cn.xlf"
$wsZhCn.Range("H5").Value = "2016-10-19 16:09:49"

foreach ($h in $wsZhCn.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$4') {
        $h.TextToDisplay = "e2c142c8-972f-4f3e-9189-484965d728c7.md"
    } elseif ($addr -eq '$A$5') {
        $h.TextToDisplay = "8de05a06-e841-430d-ad60-ba25ede17482.md"
    }
}

# ---------------------------------------------------------------------------
# "de-de" sheet (columns A:P)
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A4").Value = "e2c142c8-972f-4f3e-9189-484965d728c7.md"
$wsDeDe.Range("C4").Value = "In Translation"
$wsDeDe.Range("G4").Value = "e2c142c8-972f-4f3e-9189-484965d728c7.8bc76d304de52b61f9888771bb68e840895cbd08.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-10-19 16:09:25"

$wsDeDe.Range("A5").Value = "8de05a06-e841-430d-ad60-ba25ede17482.md"
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("G5").Value = "8de05a06-e841-430d-ad60-ba25ede17482.205adbe7658bb18bc81004be3c10cc68a87a9472.de-de.xlf"
$wsDeDe.Range("H5").Value = "2016-10-19 16:10:00"

foreach ($h in $wsDeDe.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$4') {
        $h.TextToDisplay = "e2c142c8-972f-4f3e-9189-484965d728c7.md"
    } elseif ($addr -eq '$A$5') {
        $h.TextToDisplay = "8de05a06-e841-430d-ad60-ba25ede17482.md"
    }
}
